$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'223.54"
$ws.Range("G2").Value = "'2"

$ws.Range("D3").Value = "'22.74"
$ws.Range("G3").Value = "'2"

$ws.Range("D4").Value = "'5.180"
$ws.Range("G4").Value = "'2"

$ws.Range("D5").Value = "'0.05557"
$ws.Range("G5").Value = "'2"

$ws.Range("D6").Value = "'3.389"
$ws.Range("G6").Value = "'2"

$ws.Range("D7").Value = "'6.468"
$ws.Range("G7").Value = "'2"

$ws.Range("D8").Value = "'1.089"
$ws.Range("E8").Value = "7FTXTokenFTTWorstin24h"
$ws.Range("G8").Value = "'2"

$ws.Range("D9").Value = "'0.7861"
$ws.Range("G9").Value = "'2"

$ws.Range("D10").Value = "'0.1401"
$ws.Range("G10").Value = "'2"

$ws.Range("D11").Value = "'0.07271"
$ws.Range("G11").Value = "'2"

$ws.Range("D12").Value = "'0.03143"
$ws.Range("G12").Value = "'2"

$ws.Range("D13").Value = "'0.02949"
$ws.Range("G13").Value = "'2"

$ws.Range("D14").Value = "'0.09278"
$ws.Range("G14").Value = "'2"

$ws.Range("D15").Value = "'0.001666"
$ws.Range("G15").Value = "'2"

$ws.Range("D16").Value = "'3.272"
$ws.Range("G16").Value = "'2"

$ws.Range("D17").Value = "'0.04757"
$ws.Range("G17").Value = "'2"

$ws.Range("D18").Value = "'0.0005884"
$ws.Range("G18").Value = "'2"

$ws.Range("D19").Value = "'0.006227"
$ws.Range("G19").Value = "'2"

$ws.Range("D20").Value = "'0.005262"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("G20").Value = "'2"

$ws.Range("D21").Value = "'0.001067"
$ws.Range("G21").Value = "'2"

$ws.Range("D22").Value = "'0.0001504"
$ws.Range("G22").Value = "'2"

$ws.Range("D23").Value = "'3.759"
$ws.Range("G23").Value = "'2"

$ws.Range("G24").Value = "'2"

$ws.Range("G25").Value = "'2"

$ws.Range("D26").Value = "'0.1277"
$ws.Range("G26").Value = "'2"

$ws.Range("D27").Value = "'0.0005032"
$ws.Range("E27").Value = "26UpBotsUBXT"
$ws.Range("G27").Value = "'2"

$ws.Range("G28").Value = "'2"

$ws.Range("G29").Value = "'2"

$ws.Range("G30").Value = "'2"

$ws.Range("G31").Value = "'2"

$ws.Range("G32").Value = "'2"

$ws.Range("G33").Value = "'2"

$ws.Range("G34").Value = "'2"

$ws.Range("G35").Value = "'2"

$ws.Range("G36").Value = "'2"

$ws.Range("G37").Value = "'2"

$ws.Range("G38").Value = "'2"

$ws.Range("G39").Value = "'2"

$ws.Range("D40").Value = "'0.03966"
$ws.Range("G40").Value = "'2"

$ws.Range("D41").Value = "'0.007174"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("G41").Value = "'2"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003508"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("G42").Value = "'2"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1031"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("G43").Value = "'2"

$ws.Range("D44").Value = "'0.009084"
$ws.Range("G44").Value = "'2"

$ws.Range("D45").Value = "'0.00005542"
$ws.Range("G45").Value = "'2"

$ws.Range("D46").Value = "'0.00000000752"
$ws.Range("G46").Value = "'2"

$ws.Range("D47").Value = "'0.6768"
$ws.Range("G47").Value = "'2"

$ws.Range("D48").Value = "'0.08990"
$ws.Range("G48").Value = "'2"

$ws.Range("D49").Value = "'0.00002105"
$ws.Range("G49").Value = "'2"

$ws.Range("D50").Value = "'0.01012"
$ws.Range("G50").Value = "'2"

$ws.Range("G51").Value = "'2"
